# Actualizacion Datos Personales 4 nov
# Applies the "5 de noviembre" tutor/contact-data update to the
# Cruz Alejo José Armando_2021 workbook.

$wb = $excel.ActiveWorkbook

$ws5AEM = $wb.Worksheets.Item("5AEM")
$ws5BEM = $wb.Worksheets.Item("5BEM")

# --- Sheet "5AEM" ---------------------------------------------------

# Row 7 (CARRASCO SANDOVAL, CRISTIAN ANTONIO): tutor info was missing,
# fill it in (tutor contact duplicates the student's own contact info).
$ws5AEM.Range("H7").Value = "BLANCA ESTELA SANDOVAL DÍAZ"
$ws5AEM.Range("I7").Value = "as6198099@gmail.com"
$ws5AEM.Range("J7").Value = "2721417437"

# Row 18 (PEREZ ROMERO, YAIR ANTONIO): corrected e-mail address.
$ws5AEM.Range("E18").Value = "yair26prz@gmail.com"

# Row 37 (VALDERRAMA RODRIGUEZ, EMILIO): tutor info was missing, fill it in.
$ws5AEM.Range("H37").Value = "MARÍA TEREZA RÓDRIGUEZ LOPEZ"
$ws5AEM.Range("I37").Value = "Maytequila133@gmail.com"
$ws5AEM.Range("J37").Value = "2721270249"

# --- Sheet "5BEM" ---------------------------------------------------

# Row 8 (CITLAHUA HERNANDEZ, RAUL ARTURO): tutor name corrected
# (placeholder "FINADO FINADO FINADO" replaced), and tutor phone added.
$ws5BEM.Range("H8").Value = "MARÍA FILOMENA HERNÁNDEZ CHONCOA"
$ws5BEM.Range("J8").Value = "2722480188"

# Row 13 (FLORES DE LA CRUZ, JUAN ANTONIO): corrected mobile phone
# number, and tutor name/phone added.
$ws5BEM.Range("F13").Value = "6631083766"
$ws5BEM.Range("H13").Value = "GABRIELA FLORES DE LA CRUZ"
$ws5BEM.Range("J13").Value = "2721135977"

# Row 25 (MUÑOZ LUNA, ARIAN ALEXIS): corrected tutor name.
$ws5BEM.Range("H25").Value = "JUAN MAURICIO MUÑOZ MARTINEZ"
